$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the product row with barcode 5414150631147 (row 3), shifting subsequent rows up
$ws.Rows(3).Delete()

# Update the selection to reflect the new state
$ws.Range("A3:B7").Select()
